$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The data table (date in column A, value in column B) currently ends at
# row 83. Append the next quarterly observation as row 84, copying the
# formatting (date number format, bold, border, centered) already used
# by the date column from the row above it.
$ws.Range("A83").Copy() | Out-Null
$ws.Range("A84").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A84").Value = 45884
$ws.Range("B84").Value = 0.8783323788356512
